# Edit: add "PO Forecast" sheet with forecast data, and rename the
# "Requested quantity" headers on the two existing sheets to match the
# new forecast-oriented naming scheme.

$wb = $excel.ActiveWorkbook

# --- 1) Rename header on "Weekly Quantity" sheet (sheet 1) ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Rename header on "Monthly Trend" sheet (sheet 2) ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) Add a new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (forecast output): ds, PO_Forecast, yhat_lower, yhat_upper
$forecastData = @"
45207.99999999999	112	-58.79844233487992	279.7396821760149
45214.99999999999	113	-66.62550910383793	289.6346998483162
45221.99999999999	114	-54.01821338958175	288.7549032498516
45228.99999999999	115	-54.23241335699939	290.9057263874448
45235.99999999999	116	-55.00516127442363	294.6533925954742
45242.99999999999	116	-48.59343477593214	294.1457885252713
45249.99999999999	117	-47.6056130867667	296.9280012012326
45256.99999999999	118	-70.68376412800141	304.5269824829912
45263.99999999999	119	-58.89125146481592	280.6406155491514
45270.99999999999	120	-43.82698631170107	294.8967326182299
45277.99999999999	121	-61.00019811492076	307.0970614124209
45298.99999999999	123	-53.69884646464946	299.8315119994359
45305.99999999999	124	-47.88307037047192	286.1216038722903
45312.99999999999	125	-38.20509714330692	302.2567212591962
45326.99999999999	127	-44.87287777993823	310.5777548462612
45333.99999999999	128	-44.95490314159684	304.7686186695204
45347.99999999999	130	-39.1149882580626	303.9451876694436
45361.99999999999	131	-27.86011817686937	316.5354641340138
45375.99999999999	133	-48.0904235716667	320.0475445808512
45382.99999999999	134	-39.9986497619024	307.7708479146663
45431.99999999999	140	-36.56553771649618	316.9524251069653
45445.99999999999	142	-37.32807620216333	309.4913570868071
45466.99999999999	144	-31.35752005298017	313.1171528052162
45487.99999999999	147	-23.98601711673678	321.1871161275313
45494.99999999999	148	-39.84358851608742	317.4012806555612
45501.99999999999	149	-33.09406131901742	315.7066374559143
45508.99999999999	150	-10.81350553155918	334.6344438067576
45536.99999999999	153	-22.77091098064397	326.7923069316521
45543.99999999999	154	-16.43327797235706	327.955396774263
45550.99999999999	155	-10.73845683425773	334.2774519419344
45557.99999999999	156	-9.801455158288142	319.6319497953992
45578.99999999999	158	-19.98536277838026	326.4724137488204
45585.99999999999	159	-16.47824004477787	329.2925510300843
45599.99999999999	161	-11.68643863937995	333.2770549634352
45606.99999999999	162	-22.16587975487008	334.7166232353062
45613.99999999999	163	-15.96654576471777	336.0245485283361
45620.99999999999	164	1.671824506622844	350.2769690171277
45627.99999999999	165	-11.8993544819963	331.9305246493624
45634.99999999999	165	-14.01412617920807	329.336908784934
45641.99999999999	166	-13.06982132734259	336.869148384543
45648.99999999999	167	-4.4020586550688	346.0395989849919
45655.99999999999	168	-3.07520627317512	349.2275898295051
"@

$lines = $forecastData -split "\r?\n"
$r = 2
foreach ($line in $lines) {
    if ($line.Trim() -eq "") { continue }
    $parts = $line -split "\t"
    $wsForecast.Cells.Item($r, 1).Value = [double]$parts[0]
    $wsForecast.Cells.Item($r, 2).Value = [double]$parts[1]
    $wsForecast.Cells.Item($r, 3).Value = [double]$parts[2]
    $wsForecast.Cells.Item($r, 4).Value = [double]$parts[3]
    $r++
}

$wsForecast.Range("A2:A" + ($r - 1)).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Restore original active-sheet selection (leave view state as it was
# before the new sheet was appended) instead of leaving the freshly
# added sheet selected.
[void]$wsWeekly.Select()
